$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column before column D ("inSHA_CM"), shifting old D:J -> E:K
$ws.Columns("D:D").Insert()
$ws.Range("D1").Value = "inSHA_CM"
$ws.Columns("D:D").ColumnWidth = 10

# hsGradAndBelow (row 5) is now also sourced from inSHA_CM, in addition to fromACS
$ws.Range("D5").Value = "x"

# Insert a new row before row 8 ("poverty150"), shifting old rows 8:22 -> 9:23
$ws.Rows("8:8").Insert()
$ws.Range("A8").Value = "poverty150"
$ws.Range("D8").Value = "x"
$ws.Range("E8").Value = "x"
$ws.Range("H8").Value = "neg"
$ws.Range("I8").Value = "pov"
$ws.Range("J8").Value = "Poverty"
$ws.Range("K8").Value = "Percentage of population that is Below Federal Poverty Rate of 150%"

# Update selection to match the saved state
$ws.Range("J28").Select()
